$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2024-10-05"
$ws.Range("A3").Style = "Normal"
